$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "type"
$ws.Range("C2").Value = "valdation description"
$ws.Range("C3").Value = "valdation description"
$ws.Range("C4").Value = "valdation description"
$ws.Range("C5").Value = "button"
$ws.Range("C6").Value = "link"
